# Update the "Förändrad" (Changed) date column (C) for data rows 2-32
# from serial date 45184 (2023-09-15) to 45185 (2023-09-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 32; $row++) {
    $cell = $ws.Cells.Item($row, 3)   # Column C
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45185
    }
}
